# Apply the "extraction des informations des conjectures json" edit:
#  1. Articles sheet: clear the placeholder empty DOI cells (D column) that
#     had no real DOI value, so they disappear from the sheet entirely.
#  2. Conjectures sheet: populate rows 2-23 with the per-article conjecture
#     extraction results (Article_id / Conjecture / Page).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Articles sheet — drop empty placeholder DOI cells
# ---------------------------------------------------------------------
$articles = $wb.Worksheets.Item("Articles")

$emptyDoiRows = @(2,3,4,6,7,8,9,10,11,12,14,16,17,18,19,20,21,22,23)
foreach ($r in $emptyDoiRows) {
    $articles.Cells.Item($r, 4).ClearContents()
}

# ---------------------------------------------------------------------
# 2. Conjectures sheet — fill in extracted conjecture data
# ---------------------------------------------------------------------
$conjectures = $wb.Worksheets.Item("Conjectures")

$noConjectureMissingJson = "aucune conjecture (json manquant)"
$noConjecture = "aucune conjecture"
$delta = [char]0x0394
$misConjecture = [string][char]34 + "Is there an MIS LCA with query complexity poly(" + $delta + " log n) ?" + [char]34

$data = @(
    @(1,  $noConjectureMissingJson, $null),
    @(2,  $noConjectureMissingJson, $null),
    @(3,  $noConjectureMissingJson, $null),
    @(4,  $noConjectureMissingJson, $null),
    @(5,  $noConjectureMissingJson, $null),
    @(6,  $noConjectureMissingJson, $null),
    @(7,  $noConjectureMissingJson, $null),
    @(8,  $noConjectureMissingJson, $null),
    @(9,  $noConjectureMissingJson, $null),
    @(10, $noConjectureMissingJson, $null),
    @(11, $noConjectureMissingJson, $null),
    @(12, $noConjectureMissingJson, $null),
    @(13, $noConjecture, $null),
    @(14, $noConjecture, $null),
    @(15, $noConjecture, $null),
    @(16, $noConjecture, $null),
    @(17, $noConjecture, $null),
    @(18, $noConjecture, $null),
    @(19, $noConjecture, $null),
    @(20, $misConjecture, 7),
    @(21, $noConjecture, $null),
    @(22, $noConjecture, $null)
)

$row = 2
foreach ($entry in $data) {
    $articleId   = $entry[0]
    $conjecture  = $entry[1]
    $page        = $entry[2]

    $conjectures.Cells.Item($row, 1).Value = $articleId
    $conjectures.Cells.Item($row, 2).Value = $conjecture
    if ($null -ne $page) {
        $conjectures.Cells.Item($row, 3).Value = $page
    }

    $row = $row + 1
}

# Match the reported active selection / tab on the Conjectures sheet.
$conjectures.Activate()
$conjectures.Range("E8").Select()
